$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 4660.3335
$ws.Range("I43").Value = 4660.3335
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 4660.3335
$ws.Range("L43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -4591.3335

$ws.Range("H62").Value = 2278.8572
$ws.Range("I62").Value = 2489.3333
$ws.Range("J62").Value = 1900
$ws.Range("K62").Value = 2489.3333
$ws.Range("L62").Value = 1900
$ws.Range("M62").Value = -1865.3333
$ws.Range("N62").Value = -3148

$ws.Range("H65").Value = 2278.8572
$ws.Range("I65").Value = 2489.3333
$ws.Range("J65").Value = 1900
$ws.Range("K65").Value = 12446.6665
$ws.Range("L65").Value = 9500
$ws.Range("M65").Value = -9326.666499999999
$ws.Range("N65").Value = -15740

$ws.Range("H94").Value = 6981.6665
$ws.Range("I94").Value = 3966.6667
$ws.Range("J94").Value = 9996.666999999999
$ws.Range("K94").Value = 3966.6667
$ws.Range("L94").Value = 9996.666999999999
$ws.Range("M94").Value = -3515.6667
$ws.Range("N94").Value = -10898.667

$ws.Range("H98").Value = 633.2308
$ws.Range("I98").Value = 605.55554
$ws.Range("J98").Value = 695.5
$ws.Range("K98").Value = 605.55554
$ws.Range("L98").Value = 695.5
$ws.Range("M98").Value = 892.44446
$ws.Range("N98").Value = -3691.5

$ws.Range("H112").Value = 1291.4054
$ws.Range("J112").Value = 1313.3889
$ws.Range("L112").Value = 3940.1667
$ws.Range("N112").Value = -6156.1667

$ws.Range("H113").Value = 335501.66
$ws.Range("I113").Value = 335501.66
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 335501.66
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -332247.66

$ws.Range("H122").Value = 633.2308
$ws.Range("I122").Value = 605.55554
$ws.Range("J122").Value = 695.5
$ws.Range("K122").Value = 1816.66662
$ws.Range("L122").Value = 2086.5
$ws.Range("M122").Value = 633.33338
$ws.Range("N122").Value = -6986.5

$ws.Range("H132").Value = 3911304
$ws.Range("I132").Value = 4634899.5
$ws.Range("J132").Value = 3887.9
$ws.Range("K132").Value = 13904698.5
$ws.Range("L132").Value = 11663.7
$ws.Range("M132").Value = -13902168.5
$ws.Range("N132").Value = -16723.7

$ws.Range("H141").Value = 1944.9048
$ws.Range("I141").Value = 1696.7368
$ws.Range("J141").Value = 4302.5
$ws.Range("K141").Value = 5090.2104
$ws.Range("L141").Value = 12907.5
$ws.Range("M141").Value = 89.78960000000006
$ws.Range("N141").Value = -23267.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 144849.42
$ws.Range("I2").Value = 2324.3333
$ws.Range("J2").Value = 1000000
$ws.Range("K2").Value = 2324.3333
$ws.Range("L2").Value = 1000000
$ws.Range("M2").Value = -2211.3333
$ws.Range("N2").Value = -1000226

$ws.Range("H32").Value = 32980.11
$ws.Range("I32").Value = 9405.696
$ws.Range("K32").Value = 9405.696
$ws.Range("M32").Value = -9118.696

$ws.Range("H45").Value = 94245.27
$ws.Range("I45").Value = 127742.25
$ws.Range("J45").Value = 4920
$ws.Range("K45").Value = 127742.25
$ws.Range("L45").Value = 4920
$ws.Range("M45").Value = -127365.25
$ws.Range("N45").Value = -5674

$ws.Range("H110").Value = 143157840
$ws.Range("I110").Value = 143157840
$ws.Range("K110").Value = 143157840
$ws.Range("M110").Value = -143155795

$ws.Range("H116").Value = 144849.42
$ws.Range("I116").Value = 2324.3333
$ws.Range("J116").Value = 1000000
$ws.Range("K116").Value = 2324.3333
$ws.Range("L116").Value = 1000000
$ws.Range("M116").Value = -30.33329999999978
$ws.Range("N116").Value = -1004588

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 144849.42
$ws.Range("I3").Value = 2324.3333
$ws.Range("J3").Value = 1000000
$ws.Range("K3").Value = 2324.3333
$ws.Range("L3").Value = 1000000
$ws.Range("M3").Value = -2210.3333
$ws.Range("N3").Value = -1000228

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 63238.25
$ws.Range("I16").Value = 738.4167
$ws.Range("J16").Value = 250737.75
$ws.Range("K16").Value = 738.4167
$ws.Range("L16").Value = 250737.75
$ws.Range("M16").Value = -451.4167
$ws.Range("N16").Value = -251311.75

$ws.Range("H99").Value = 7070.304
$ws.Range("I99").Value = 2298
$ws.Range("J99").Value = 8075
$ws.Range("K99").Value = 2298
$ws.Range("L99").Value = 8075
$ws.Range("M99").Value = -800
$ws.Range("N99").Value = -11071

$ws.Range("H111").Value = 40000
$ws.Range("J111").Value = 40000
$ws.Range("L111").Value = 40000
$ws.Range("N111").Value = -48180

$ws.Range("H113").Value = 63238.25
$ws.Range("I113").Value = 738.4167
$ws.Range("J113").Value = 250737.75
$ws.Range("K113").Value = 738.4167
$ws.Range("L113").Value = 250737.75
$ws.Range("M113").Value = 1431.5833
$ws.Range("N113").Value = -255077.75

$ws.Range("H122").Value = 1191.4166
$ws.Range("I122").Value = 1190.7273
$ws.Range("J122").Value = 1199
$ws.Range("K122").Value = 3572.1819
$ws.Range("L122").Value = 3597
$ws.Range("M122").Value = -1122.1819
$ws.Range("N122").Value = -8497

$ws.Range("H126").Value = 7070.304
$ws.Range("I126").Value = 2298
$ws.Range("J126").Value = 8075
$ws.Range("K126").Value = 6894
$ws.Range("L126").Value = 24225
$ws.Range("M126").Value = -4424
$ws.Range("N126").Value = -29165

$ws.Range("H134").Value = 1394.8572
$ws.Range("I134").Value = 553.38464
$ws.Range("J134").Value = 2762.25
$ws.Range("K134").Value = 1660.15392
$ws.Range("L134").Value = 8286.75
$ws.Range("M134").Value = 874.84608
$ws.Range("N134").Value = -13356.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 6572.5264
$ws.Range("I80").Value = 966.6667
$ws.Range("J80").Value = 7623.625
$ws.Range("K80").Value = 2900.0001
$ws.Range("L80").Value = 22870.875
$ws.Range("M80").Value = -1964.0001
$ws.Range("N80").Value = -24742.875

$ws.Range("H83").Value = 6572.5264
$ws.Range("I83").Value = 966.6667
$ws.Range("J83").Value = 7623.625
$ws.Range("K83").Value = 8700.0003
$ws.Range("L83").Value = 68612.625
$ws.Range("M83").Value = -4020.0003
$ws.Range("N83").Value = -77972.625

$ws.Range("H99").Value = 1467.1428
$ws.Range("I99").Value = 878.3333
$ws.Range("J99").Value = 5000
$ws.Range("K99").Value = 2634.9999
$ws.Range("L99").Value = 15000
$ws.Range("M99").Value = -388.9998999999998
$ws.Range("N99").Value = -19492

$ws.Range("H107").Value = 861.5526
$ws.Range("I107").Value = 907.5
$ws.Range("J107").Value = 852.9375
$ws.Range("K107").Value = 2722.5
$ws.Range("L107").Value = 2558.8125
$ws.Range("M107").Value = -802.5
$ws.Range("N107").Value = -6398.8125

$ws.Range("H108").Value = 1616.3334
$ws.Range("I108").Value = 1217.4
$ws.Range("K108").Value = 3652.2
$ws.Range("M108").Value = -772.2000000000003

$ws.Range("H120").Value = 1000030
$ws.Range("I120").Value = 1000030
$ws.Range("K120").Value = 3000090
$ws.Range("M120").Value = -2995252

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 500
$ws.Range("I113").Value = 500
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 500
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = 1670

$ws.Range("H132").Value = 3361.0645
$ws.Range("I132").Value = 2540.318
$ws.Range("K132").Value = 7620.954000000001
$ws.Range("M132").Value = -5090.954000000001

$ws.Range("H136").Value = 395953.75
$ws.Range("J136").Value = 395953.75
$ws.Range("L136").Value = 1187861.25
$ws.Range("N136").Value = -1192961.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2279.5
$ws.Range("J100").Value = 2419.25
$ws.Range("L100").Value = 2419.25
$ws.Range("N100").Value = -3501.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 62501224
$ws.Range("I96").Value = 142858960
$ws.Range("J96").Value = 766.1111
$ws.Range("K96").Value = 142858960
$ws.Range("L96").Value = 766.1111
$ws.Range("M96").Value = -142857587
$ws.Range("N96").Value = -3512.1111

$ws.Range("H122").Value = 2735.6843
$ws.Range("I122").Value = 1998.2
$ws.Range("J122").Value = 5501.25
$ws.Range("K122").Value = 5994.6
$ws.Range("L122").Value = 16503.75
$ws.Range("M122").Value = -3544.6
$ws.Range("N122").Value = -21403.75

$ws.Range("H137").Value = 47000
$ws.Range("J137").Value = 47000
$ws.Range("L137").Value = 47000
$ws.Range("N137").Value = -57200
